$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.448244571685791
$ws.Range("B1").Value = 3.487921237945557
$ws.Range("C1").Value = 2.740448236465454
$ws.Range("D1").Value = 2.265465974807739
$ws.Range("E1").Value = 1.609159350395203
